# ---------------------------------------------------------------------------
# Applies the "Mise à jour de la carte avec nouveaux sites" edit:
#   - adds centered alignment (style index 1) to the PS_latitude / PS_longitude
#     columns (J:K) for the existing data rows 5-120 (including the blank
#     separator rows, which gain empty styled J/K cells)
#   - appends 8 new site rows (121-128) with Site / Latitude / Longitude /
#     PS_plus_proche / PS_latitude / PS_longitude / Distance_km data
#   - updates the sheet selection
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108
$xlPasteFormats = -4122

# Helper: write a value as a genuine *text* cell (shared string), even when
# the text looks like a number (Excel would otherwise silently convert
# "43.746..." into a numeric value and lose the original text formatting).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value2 = $value
}

# ---------------------------------------------------------------------------
# 1. Style the PS_latitude (J) / PS_longitude (K) columns for rows 5-120
#    (center horizontal + vertical, matching style index 1 already used by
#    the other columns on these rows). A single range covers both the rows
#    that already hold data and the blank separator rows (12, 80, 112, 113,
#    115, 118, 119, 120), which simply gain empty, styled J/K cells.
# ---------------------------------------------------------------------------
$jk = $ws.Range("J5:K120")
$jk.HorizontalAlignment = $xlCenter
$jk.VerticalAlignment = $xlCenter

# ---------------------------------------------------------------------------
# 2. Append the new site rows (121-128)
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row = 121; F = "soustons";            G = "43.74641770860991";  H = " -1.33285379210912";   I = "Soustons";             J = "43.771091126263514";  K = " -1.2948007929781047"; L = "4.111690" },
    @{ Row = 122; F = "mimizan";              G = "44.19725528555141";  H = "-1.2291033861765632";   I = "Mimizan";              J = "44.1948017922608";    K = " -1.2295049995137841"; L = "0.274506" },
    @{ Row = 123; F = "morcenx la nouvelle";  G = "44.02969680821508";  H = "-0.9227099943452951";   I = "Morcenx-la-Nouvelle";  J = "44.03143778431751";   K = " -0.9052865619625207"; L = "1.410082" },
    @{ Row = 124; F = "lit et mixe";          G = "44.02549860293042";  H = "-1.2583691579531584";   I = "Linxe";                J = "43.907624274812235";  K = " -1.2321315191950484"; L = "13.265380" },
    @{ Row = 125; F = "saint sever";          G = "43.7502878518836";   H = "-0.568674355747199";    I = "Saint-Sever";          J = "43.76397390663459";   K = " -0.5258372653500577"; L = "3.770054" },
    @{ Row = 126; F = "dax";                  G = "43.70987517571523";  H = "-1.0441713878412118";   I = "Saint-Paul-lès-Dax";   J = "43.7263292279978";    K = " -1.03918106238464";   L = "1.871864" },
    @{ Row = 127; F = "PEYREHORADE";          G = "43.5512456094283";   H = "-1.1301652172450665";   I = "Guiche";               J = "43.512432834195224";  K = " -1.2207089798077064"; L = "8.494906" },
    @{ Row = 128; F = "PEYREHORADE soumo";    G = "43.55076012419594";  H = "-1.1299851875293643";   I = "Guiche";               J = "43.512432834195224";  K = " -1.2207089798077064"; L = "8.480260" }
)

# Use L5 (style index 1: centered horizontal + vertical) as the format donor
# for the F, G, H, I and L cells of the new rows.
$ws.Range("L5").Copy()

foreach ($r in $newRows) {
    $row = $r.Row

    foreach ($col in @("F", "G", "H", "I", "L")) {
        $addr = $col + $row
        $rng = $ws.Range($addr)
        Set-TextValue $rng $r[$col]
        $rng.PasteSpecial($xlPasteFormats)
    }

    # J (PS_latitude) / K (PS_longitude) keep the default (unstyled) look,
    # exactly like the other rows before this edit.
    foreach ($col in @("J", "K")) {
        $addr = $col + $row
        $rng = $ws.Range($addr)
        Set-TextValue $rng $r[$col]
        $rng.Style = "Normal"
    }
}

# ---------------------------------------------------------------------------
# 3. Update the view / selection
# ---------------------------------------------------------------------------
$ws.Range("L120").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 103
$win.ScrollColumn = 8
